$wb = $excel.ActiveWorkbook

$wsBinary = $wb.Worksheets.Item("factor_binary")
$wsFactor3 = $wb.Worksheets.Item("factor_3")

# New row 23
$wsBinary.Range("B23").Value = "false (correct)"
$wsBinary.Range("C23").Value = "true."
$wsBinary.Range("A23").Value = "true2"

# New row 24
$wsBinary.Range("B24").Value = "false."
$wsBinary.Range("C24").Value = "true (correct)"
$wsBinary.Range("A24").Value = "true3"

# B23 picks up the "Lucida Console / vertical-center" style already used
# elsewhere in the workbook (factor_3!C12) instead of Excel fabricating a
# brand-new cellXf.
$wsFactor3.Range("C12").Copy()
$wsBinary.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active sheet/selection from factor_3 back to factor_binary.
$wsBinary.Activate()
$wsBinary.Application.ActiveWindow.ScrollRow = 13
$wsBinary.Range("A25").Select()
